$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (M2:T2)
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 5.769525850532444
$ws.Range("R2").Value = 51.925732654792
$ws.Range("S2").Value = 0.2299953477621856
$ws.Range("T2").Value = 0.2299953477621856

# Row 3 updates (only O3, P3, S3, T3 change)
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("S3").Value = 0.6794731949692173
$ws.Range("T3").Value = 0.6794731949692174

# Row 4 updates (M4:T4)
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 2.271018036145778
$ws.Range("R4").Value = 20.439162325312
$ws.Range("S4").Value = 0.09053145726859702
$ws.Range("T4").Value = 0.09053145726859703
